# Scheduled runner update: refresh market-board derived values (current
# average prices, computed leve costs/profits) across the job sheets.
#
# For each worksheet, a small set of rows gets new values in columns
# H (currentAveragePrice), I (currentAveragePriceNQ), J (currentAveragePriceHQ),
# K (LevePriceNQ), L (LevePriceHQ), M (LeveProfitNQ), N (LeveProfitHQ).

$wb = $excel.ActiveWorkbook

# NOTE: this interpreter's function dispatch only resolves *positional*
# parameters reliably, so Set-RowValues takes plain positional args
# (sheet name, row number, hashtable of column -> new value) instead of
# PowerShell named parameters.
function Set-RowValues($SheetName, $Row, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# ---- ALC ----
Set-RowValues "ALC" 33 @{
    H = 214.78572; I = 100.71429; J = 328.85715; K = 100.71429; L = 328.85715; M = 128.28571; N = -786.85715
}
Set-RowValues "ALC" 64 @{
    H = 3932.3794; I = 3761.7932; K = 3761.7932; M = -3513.7932
}
Set-RowValues "ALC" 67 @{
    H = 3932.3794; I = 3761.7932; K = 3761.7932; M = -2903.7932
}
Set-RowValues "ALC" 76 @{
    H = 3532.889; I = 3500; J = 3534.1538; K = 3500; L = 3534.1538; M = -3185; N = -4164.1538
}
Set-RowValues "ALC" 79 @{
    H = 3532.889; I = 3500; J = 3534.1538; K = 3500; L = 3534.1538; M = -2408; N = -5718.1538
}
Set-RowValues "ALC" 112 @{
    H = 1310.8; I = 700; J = 1389.6129; K = 2100; L = 4168.8387; M = -992; N = -6384.8387
}
Set-RowValues "ALC" 129 @{
    H = 709.46; I = 465.66666; J = 717; K = 1396.99998; L = 2151; M = 3603.00002; N = -12151
}
Set-RowValues "ALC" 138 @{
    H = 2835.303; I = 1144.3636; J = 3046.6704; K = 3433.0908; L = 9140.0112; M = 1706.9092; N = -19420.0112
}

# ---- ARM ----
Set-RowValues "ARM" 6 @{
    H = 7500; I = 5000; J = 10000; K = 5000; L = 10000; M = -4827; N = -10346
}
Set-RowValues "ARM" 63 @{
    H = 2435.3333; I = 1370.5714; J = 3367; K = 1370.5714; L = 3367; M = -684.5714; N = -4739
}
Set-RowValues "ARM" 66 @{
    H = 2435.3333; I = 1370.5714; J = 3367; K = 6852.857; L = 16835; M = -3420.857; N = -23699
}

# ---- BSM ----
Set-RowValues "BSM" 54 @{
    H = 611; I = 611; K = 611; M = -127
}
Set-RowValues "BSM" 105 @{
    H = 2327.6667; I = 1816.6666; J = 2838.6667; K = 1816.6666; L = 2838.6667; M = -69.6666; N = -6332.6667
}

# ---- CRP ----
Set-RowValues "CRP" 16 @{
    H = 2600; I = 1900; J = 3066.6667; K = 1900; L = 3066.6667; M = -1613; N = -3640.6667
}
Set-RowValues "CRP" 31 @{
    H = 744174.4; I = 2116.9092; J = 1910264.6; K = 2116.9092; L = 1910264.6; M = -1821.9092; N = -1910854.6
}
Set-RowValues "CRP" 34 @{
    H = 744174.4; I = 2116.9092; J = 1910264.6; K = 2116.9092; L = 1910264.6; M = -1914.9092; N = -1910668.6
}
Set-RowValues "CRP" 62 @{
    H = 2962.5; I = 2520; J = 3700; K = 2520; L = 3700; M = -1896; N = -4948
}
Set-RowValues "CRP" 65 @{
    H = 2962.5; I = 2520; J = 3700; K = 12600; L = 18500; M = -9480; N = -24740
}
Set-RowValues "CRP" 113 @{
    H = 2600; I = 1900; J = 3066.6667; K = 1900; L = 3066.6667; M = 270; N = -7406.6667
}
Set-RowValues "CRP" 122 @{
    H = 2099.625; I = 1599; J = 2400; K = 4797; L = 7200; M = -2347; N = -12100
}

# ---- CUL ----
Set-RowValues "CUL" 122 @{
    H = 682236.7; I = 4967.8335; J = 2393231.8; K = 44710.5015; L = 21539086.2; M = -42260.5015; N = -21543986.2
}

# ---- GSM ----
Set-RowValues "GSM" 10 @{
    H = 500015000; I = 500015000; K = 500015000; M = -500014831
}
Set-RowValues "GSM" 70 @{
    H = 6902.6313; I = 7063.8887; J = 4000; K = 7063.8887; L = 4000; M = -6793.8887; N = -4540
}
Set-RowValues "GSM" 73 @{
    H = 6902.6313; I = 7063.8887; J = 4000; K = 7063.8887; L = 4000; M = -6127.8887; N = -5872
}
Set-RowValues "GSM" 80 @{
    H = 3015.413; I = 2815.6316; J = 3156; K = 2815.6316; L = 3156; M = -1817.6316; N = -5152
}
Set-RowValues "GSM" 83 @{
    H = 3015.413; I = 2815.6316; J = 3156; K = 14078.158; L = 15780; M = -9086.158; N = -25764
}
Set-RowValues "GSM" 100 @{
    H = 17354; J = 17354; L = 17354; N = -19518
}
Set-RowValues "GSM" 113 @{
    H = 720460.5; I = 1400; K = 1400; M = 770
}

# ---- LTW ----
Set-RowValues "LTW" 68 @{
    H = 2013.931; J = 1972.7273; L = 1972.7273; N = -3470.7273
}
Set-RowValues "LTW" 71 @{
    H = 2013.931; J = 1972.7273; L = 9863.6365; N = -17351.6365
}

# ---- WVR ----
Set-RowValues "WVR" 62 @{
    H = 11136462; I = 20041540; J = 5113.25; K = 20041540; L = 5113.25; M = -20040916; N = -6361.25
}
Set-RowValues "WVR" 65 @{
    H = 11136462; I = 20041540; J = 5113.25; K = 100207700; L = 25566.25; M = -100204580; N = -31806.25
}
Set-RowValues "WVR" 99 @{
    H = 45048.668; J = 45048.668; L = 45048.668; N = -51038.668
}
Set-RowValues "WVR" 100 @{
    H = 689.38464; I = 466.8889; J = 1190; K = 933.7778; L = 2380; M = -392.7778; N = -3462
}
Set-RowValues "WVR" 122 @{
    H = 3951.4; I = 4918.8; J = 2984; K = 14756.4; L = 8952; M = -12306.4; N = -13852
}
Set-RowValues "WVR" 126 @{
    H = 2667.0938; I = 2589.5; J = 2899.875; K = 7768.5; L = 8699.625; M = -5298.5; N = -13639.625
}
